$p = $ppt.ActivePresentation

# --- 1) Update cached date-field text ("9/12/2017" -> "12/21/2018") ---
#     Slide Master date placeholder
$p.SlideMaster.Shapes.Item(2).TextFrame.TextRange.Text = "12/21/2018"

#     Slide Layout 1 (Title Slide) date placeholder
$p.SlideMaster.CustomLayouts.Item(1).Shapes.Item(3).TextFrame.TextRange.Text = "12/21/2018"

#     Slide Layout 2 date placeholder
$p.SlideMaster.CustomLayouts.Item(2).Shapes.Item(2).TextFrame.TextRange.Text = "12/21/2018"

#     Slide Layout 3 date placeholder
$p.SlideMaster.CustomLayouts.Item(3).Shapes.Item(3).TextFrame.TextRange.Text = "12/21/2018"

#     Notes Master date placeholder
$p.NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text = "12/21/2018"

# --- 2) Fix wording on slide 10 ("Notes") : APPLICATIONS default value EACH -> ALL ---
$slide10 = $p.Slides.Item(10)
$notesShape = $slide10.Shapes.Item(2)
$appParagraph = $notesShape.TextFrame.TextRange.Paragraphs(2)
$descRun = $appParagraph.Runs(2)
$descRun.Text = ": if no information filled, then default value "
$descRun.InsertAfter("is “ALL`"")
